$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.843.50"

$ws.Range("D3").Value = "1.894.19"
$ws.Range("E3").Value = "  -0.92%  "

$ws.Range("D4").Value = "'0.9991"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").Value = "'0.7658"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.41%  "

$ws.Range("D6").Value = "'239.30"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.89%  "

$ws.Range("D7").Value = "'0.9998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("D8").Value = "'0.3046"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.22%  "

$ws.Range("D9").Value = "'25.37"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -6.32%  "

$ws.Range("D10").Value = "'0.06816"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.28%  "

$ws.Range("D11").Value = "'0.07990"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.26%  "

$ws.Range("D12").Value = "'0.7488"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.53%  "

$ws.Range("D13").Value = "1.885.28"
$ws.Range("E13").Value = "  -1.79%  "

$ws.Range("D14").Value = "'5.201"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.55%  "

$ws.Range("D15").Value = "'91.09"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.74%  "

$ws.Range("D16").Value = "29.836.38"
$ws.Range("E16").Value = "  -1.44%  "

$ws.Range("B17").Value = "Avalanche"
$ws.Range("C17").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D17").Value = "'13.87"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.98%  "

$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").Value = "'5.966"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.33%  "

$ws.Range("D19").Value = "'0.000007670"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.02%  "

$ws.Range("D20").Value = "'234.79"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.19%  "

$ws.Range("D21").Value = "'0.9999"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.06%  "

$ws.Range("D22").Value = "2.133.24"
$ws.Range("E22").Value = "  -2.39%  "

$ws.Range("D23").Value = "'0.9998"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.02%  "

$ws.Range("D24").Value = "'6.932"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.17%  "

$ws.Range("D25").Value = "'9.246"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.93%  "

$ws.Range("D26").Value = "'165.39"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.16%  "

$ws.Range("D27").Value = "'18.70"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.84%  "

$ws.Range("D28").Value = "'0.1299"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.28%  "

$ws.Range("D29").Value = "'2.043"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.26%  "

$ws.Range("D30").Value = "'1.342"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.53%  "

$ws.Range("E31").Value = "  -2.34%  "

$ws.Range("D32").Value = "'4.276"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.07%  "

$ws.Range("D33").Value = "'4.024"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.48%  "

$ws.Range("D34").Value = "'0.05364"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.69%  "

$ws.Range("D35").Value = "'1.249"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.87%  "

$ws.Range("D36").Value = "'0.7271"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.45%  "

$ws.Range("D37").Value = "'2.712"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.67%  "

$ws.Range("D38").Value = "'0.01924"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.03%  "

$ws.Range("E39").Value = "  -0.55%  "

$ws.Range("D40").Value = "'6.183"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.48%  "

$ws.Range("D41").Value = "'0.4404"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.18%  "

$ws.Range("D42").Value = "'72.23"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.02%  "

$ws.Range("E43").Value = "  -2.01%  "

$ws.Range("D44").Value = "'1.0000"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.08%  "

$ws.Range("D45").Value = "'0.8239"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.95%  "

$ws.Range("D46").Value = "'101.06"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.38%  "

$ws.Range("D47").Value = "'7.586"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.71%  "

$ws.Range("D48").Value = "'9.770"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.32%  "

$ws.Range("D49").Value = "2.035.69"
$ws.Range("E49").Value = "  -2.33%  "

$ws.Range("D50").Value = "'36.14"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.62%  "

$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "'922.10"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.07%  "
